# Bug #25 ("Drop Section") added to the "Bug Log" sheet, mirroring the
# existing "Drop Bid" row (row 26) one row below it, plus a couple of
# view/selection tweaks that went along with the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log")
$ws.Activate()

# --- New row 27: "Drop Section" bug entry ------------------------------
# Row 26 (bug #24, "Drop Bid") already carries the exact formatting
# (borders / fills / number format / wrap) that the new row needs, so
# copy its formats down one row before filling in the new values.
$srcRow = $ws.Range("A26:H26")
$dstRow = $ws.Range("A27:H27")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = 3
$ws.Cells.Item(27, 3).Value = "Drop Section"
$ws.Cells.Item(27, 4).Value = "You can drop a section once the window is closed. You should only be able to drop sections during active rounds"
$ws.Cells.Item(27, 5).Value = "Resolved"
$ws.Cells.Item(27, 6).Value = 43780
$ws.Cells.Item(27, 7).Value = 43780
$ws.Cells.Item(27, 8).Value = "Matthew & Sheng Qin"

# Match row 26's height (30.75pt / thick-bottom look) on the new row.
$ws.Rows(27).RowHeight = 30.75

# --- View / selection tweaks --------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D38").Select()
